{"js": "// Update the 25 populated cells of the single 5-column practice table\n// (3-digit \u00f7 1-digit problems) to the new set of problems. Edits are\n// addressed by (rowIndex, colIndex) so they are unambiguous even where\n// the new text duplicates a value used elsewhere in the table.\nconst edits = [\n  [0, 0, \"395\u00f73=\"], [0, 1, \"749\u00f73=\"], [0, 2, \"336\u00f75=\"], [0, 3, \"685\u00f79=\"], [0, 4, \"405\u00f75=\"],\n  [4, 0, \"583\u00f75=\"], [4, 1, \"444\u00f73=\"], [4, 2, \"747\u00f74=\"], [4, 3, \"695\u00f74=\"], [4, 4, \"254\u00f79=\"],\n  [8, 0, \"585\u00f74=\"], [8, 1, \"484\u00f73=\"], [8, 2, \"292\u00f72=\"], [8, 3, \"675\u00f72=\"], [8, 4, \"125\u00f74=\"],\n  [12, 0, \"902\u00f74=\"], [12, 1, \"877\u00f72=\"], [12, 2, \"137\u00f79=\"], [12, 3, \"685\u00f79=\"], [12, 4, \"966\u00f74=\"],\n  [16, 0, \"292\u00f77=\"], [16, 1, \"140\u00f79=\"], [16, 2, \"471\u00f78=\"], [16, 3, \"619\u00f75=\"], [16, 4, \"235\u00f78=\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nfor (const [row, col, newText] of edits) {\n  const cell = table.getCell(row, col);\n  const paragraphs = cell.body.paragraphs;\n  paragraphs.load(\"items\");\n  await context.sync();\n\n  const paragraph = paragraphs.items[0];\n  // Replace just the paragraph's text range so the existing run/paragraph\n  // formatting (font, size, alignment) carried by the original text is kept.\n  const range = paragraph.getRange();\n  range.insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Update the 25 populated cells of the single 5-column practice table\n# (3-digit \u00f7 1-digit problems) to the new set of problems. Edits are\n# addressed by (row, column) -- 1-based, as Word COM expects -- so they\n# are unambiguous even where the new text duplicates a value used\n# elsewhere in the table. Setting Cell.Range.Text keeps the existing\n# run/paragraph formatting (font, size, alignment) already on that cell.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$edits = @(\n  @(1, 1, \"395\u00f73=\"), @(1, 2, \"749\u00f73=\"), @(1, 3, \"336\u00f75=\"), @(1, 4, \"685\u00f79=\"), @(1, 5, \"405\u00f75=\"),\n  @(5, 1, \"583\u00f75=\"), @(5, 2, \"444\u00f73=\"), @(5, 3, \"747\u00f74=\"), @(5, 4, \"695\u00f74=\"), @(5, 5, \"254\u00f79=\"),\n  @(9, 1, \"585\u00f74=\"), @(9, 2, \"484\u00f73=\"), @(9, 3, \"292\u00f72=\"), @(9, 4, \"675\u00f72=\"), @(9, 5, \"125\u00f74=\"),\n  @(13, 1, \"902\u00f74=\"), @(13, 2, \"877\u00f72=\"), @(13, 3, \"137\u00f79=\"), @(13, 4, \"685\u00f79=\"), @(13, 5, \"966\u00f74=\"),\n  @(17, 1, \"292\u00f77=\"), @(17, 2, \"140\u00f79=\"), @(17, 3, \"471\u00f78=\"), @(17, 4, \"619\u00f75=\"), @(17, 5, \"235\u00f78=\")\n)\n\nforeach ($edit in $edits) {\n    $row = $edit[0]\n    $col = $edit[1]\n    $newText = $edit[2]\n    $t.Cell($row, $col).Range.Text = $newText\n}\n"}
